# Generate Report for Handoff
# Replaces the two tracked files (3d898684...md, 56e3ed2d...md) with a new
# set of three tracked files: a .md (9242a8f1...md) plus two .png deps
# (c2149f5a...png, f969858f...png), across the Overview / zh-cn / de-de
# sheets, and appends a new row 4 for the third tracked file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Drop existing hyperlinks so we can rebuild them in the exact order
# the new state needs (Hyperlinks.Add appends, and rId/order follows
# insertion order).
$ov.Hyperlinks.Delete()

$ov.Range("A2:D4").Clear()

$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-25-17 03:25:37"

$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-25-17 03:25:37"

$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-25-17 03:25:37"

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/9242a8f1-9c03-485c-8208-68ec8e196dbc.md", "", "", "9242a8f1-9c03-485c-8208-68ec8e196dbc.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/c2149f5a-fd78-40eb-b909-8d0fc080cac5.png", "", "", "c2149f5a-fd78-40eb-b909-8d0fc080cac5.png")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/f969858f-1897-451b-8867-ab993ade7621.png", "", "", "f969858f-1897-451b-8867-ab993ade7621.png")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Delete()
$zh.Range("A2:K4").Clear()

# Row 2 - the .md file, handoff still "Include"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("E2").Value = "2016-03-17 03:25:29"
$zh.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H2").Value = "0001-01-01 00:00:00"
$zh.Range("I2").Value = "Include"

# Row 3 - the first .png dependency
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("E3").Value = "2016-03-17 03:25:29"
$zh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("I3").Value = "IsDependency"
$zh.Range("J3").Value = "e2e\9242a8f1-9c03-485c-8208-68ec8e196dbc.md"

# Row 4 - the second .png dependency
$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("E4").Value = "2016-03-17 03:25:29"
$zh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H4").Value = "0001-01-01 00:00:00"
$zh.Range("I4").Value = "IsDependency"
$zh.Range("J4").Value = "e2e\9242a8f1-9c03-485c-8208-68ec8e196dbc.md"

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/9242a8f1-9c03-485c-8208-68ec8e196dbc.md", "", "", "9242a8f1-9c03-485c-8208-68ec8e196dbc.md")
$zh.Hyperlinks.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/9242a8f1-9c03-485c-8208-68ec8e196dbc.md", "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e3a1c113d8d2aef05bd6809a2d0157ed7063af5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/9242a8f1-9c03-485c-8208-68ec8e196dbc.0552ef66659e040964c1f3893adc49530494be08.zh-cn.xlf", "", "", "9242a8f1-9c03-485c-8208-68ec8e196dbc.0552ef66659e040964c1f3893adc49530494be08.zh-cn.xlf")

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/c2149f5a-fd78-40eb-b909-8d0fc080cac5.png", "", "", "c2149f5a-fd78-40eb-b909-8d0fc080cac5.png")
$zh.Hyperlinks.Add($zh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/c2149f5a-fd78-40eb-b909-8d0fc080cac5.png", "", "", ".png")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e3a1c113d8d2aef05bd6809a2d0157ed7063af5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/1e431fc6e437bd41a01bca6bd0dfd6e108685851.png", "", "", "1e431fc6e437bd41a01bca6bd0dfd6e108685851.png")

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/f969858f-1897-451b-8867-ab993ade7621.png", "", "", "f969858f-1897-451b-8867-ab993ade7621.png")
$zh.Hyperlinks.Add($zh.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/f969858f-1897-451b-8867-ab993ade7621.png", "", "", ".png")
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e3a1c113d8d2aef05bd6809a2d0157ed7063af5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/493d1825cd95f130d642e1af488986ff0c0d7680.png", "", "", "493d1825cd95f130d642e1af488986ff0c0d7680.png")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Delete()
$de.Range("A2:K4").Clear()

# Row 2 - the .md file, handoff still "Include"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("E2").Value = "2016-03-17 03:25:37"
$de.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H2").Value = "0001-01-01 00:00:00"
$de.Range("I2").Value = "Include"

# Row 3 - the first .png dependency
$de.Range("C3").Value = "Ready for handoff"
$de.Range("E3").Value = "2016-03-17 03:25:37"
$de.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("I3").Value = "IsDependency"
$de.Range("J3").Value = "e2e\9242a8f1-9c03-485c-8208-68ec8e196dbc.md"

# Row 4 - the second .png dependency
$de.Range("C4").Value = "Ready for handoff"
$de.Range("E4").Value = "2016-03-17 03:25:37"
$de.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H4").Value = "0001-01-01 00:00:00"
$de.Range("I4").Value = "IsDependency"
$de.Range("J4").Value = "e2e\9242a8f1-9c03-485c-8208-68ec8e196dbc.md"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/9242a8f1-9c03-485c-8208-68ec8e196dbc.md", "", "", "9242a8f1-9c03-485c-8208-68ec8e196dbc.md")
$de.Hyperlinks.Add($de.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/9242a8f1-9c03-485c-8208-68ec8e196dbc.md", "", "", ".md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/65d07f5e9a779efb85ea7850fe0f289c136c2e69/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/9242a8f1-9c03-485c-8208-68ec8e196dbc.0552ef66659e040964c1f3893adc49530494be08.de-de.xlf", "", "", "9242a8f1-9c03-485c-8208-68ec8e196dbc.0552ef66659e040964c1f3893adc49530494be08.de-de.xlf")

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/c2149f5a-fd78-40eb-b909-8d0fc080cac5.png", "", "", "c2149f5a-fd78-40eb-b909-8d0fc080cac5.png")
$de.Hyperlinks.Add($de.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/c2149f5a-fd78-40eb-b909-8d0fc080cac5.png", "", "", ".png")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/65d07f5e9a779efb85ea7850fe0f289c136c2e69/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/1e431fc6e437bd41a01bca6bd0dfd6e108685851.png", "", "", "1e431fc6e437bd41a01bca6bd0dfd6e108685851.png")

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/f969858f-1897-451b-8867-ab993ade7621.png", "", "", "f969858f-1897-451b-8867-ab993ade7621.png")
$de.Hyperlinks.Add($de.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/f969858f-1897-451b-8867-ab993ade7621.png", "", "", ".png")
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/65d07f5e9a779efb85ea7850fe0f289c136c2e69/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/493d1825cd95f130d642e1af488986ff0c0d7680.png", "", "", "493d1825cd95f130d642e1af488986ff0c0d7680.png")
